# "stop the dedicated offset codes. we can encode this into the operations"
#
# The dedicated SOFF*/ROFF*/OSOFF* offset opcode rows (81-89) on the
# "Opcodes" sheet are no longer needed, so their contents are cleared
# (values, string references and formulas alike), while leaving the
# existing cell styling/formatting untouched - this mirrors how row 90
# (already blank) looks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opcodes")
$ws.Activate()

# Clear out the now-unused opcode rows (79/80 numeric rows correspond to
# decimal opcodes 79-87, i.e. spreadsheet rows 81-89).
$ws.Range("A81:H89").ClearContents()

# Restore the cursor/selection to where the author left it after the edit.
$ws.Range("A81").Select()
